# Update refreshed price/volume figures on the crypto tracker sheet.
# Source values are plain text (Price/Volume columns hold strings such as
# "330.35" or "1.26%"), so force the cell format to Text before writing
# the new value - otherwise Excel would coerce them into numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "330.35";     "E2"  = "1.26%"
    "D3"  = "44.25";      "E3"  = "0.40%"
    "D4"  = "5.477";      "E4"  = "-2.13%"
    "D5"  = "0.08028";    "E5"  = "-0.19%"
    "D6"  = "2.000";      "E6"  = "5.82%"
    "D7"  = "0.9524";     "E7"  = "0.83%"
    "D8"  = "2.559";      "E8"  = "-2.73%"
                           "E9"  = "-1.62%"
    "D10" = "0.1900";     "E10" = "3.11%"
                           "E11" = "29.08%"
    "D12" = "0.09991";    "E12" = "2.24%"
    "D13" = "0.04831";    "E13" = "13.11%"
    "D14" = "0.1064";     "E14" = "-0.24%"
    "D15" = "0.001270";   "E15" = "-0.16%"
    "D16" = "0.04075";    "E16" = "-3.40%"
    "D17" = "0.005919";   "E17" = "-1.00%"
    "D18" = "3.366";      "E18" = "-6.64%"
    "D19" = "4.394";      "E19" = "2.14%"
    "D20" = "0.3430";     "E20" = "-1.88%"
                           "E21" = "1.31%"
    "D22" = "0.2502";     "E22" = "-5.81%"
                           "E23" = "2.10%"
    "D24" = "0.004368";   "E24" = "-3.17%"
                           "E25" = "-4.93%"
    "D26" = "0.0003740";  "E26" = "-6.45%"
    "D38" = "0.02600";    "E38" = "-0.88%"
    "D39" = "0.05792";    "E39" = "6.36%"
    "D40" = "0.007575";   "E40" = "-0.17%"
                           "E41" = "0.73%"
    "D42" = "0.007305";   "E42" = "-0.74%"
                           "E43" = "-1.02%"
    "D44" = "0.008829";   "E44" = "-0.07%"
    "D45" = "0.00006976"; "E45" = "0.69%"
                           "E46" = "-0.32%"
    "D47" = "0.0005795";  "E47" = "-0.29%"
                           "E48" = "54.99%"
                           "E49" = "-6.10%"
                           "E50" = "-0.32%"
                           "E51" = "-0.32%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
